# Mexico Liga MX Femenil - base update (15-04-2024 22:35)
#
# 1) A number of existing fixture rows had their betting-data columns
#    (B, F:AC) swapped between two adjacent rows (the row "id" in column A
#    plus C/D/E stay put - only the match record moved).
# 2) Five brand-new fixture rows were appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param([int]$r1, [int]$r2)

    $tmpB   = $ws.Range("B$r1").Value2
    $tmpFAC = $ws.Range("F$r1`:AC$r1").Value2

    $ws.Range("B$r1").Value        = $ws.Range("B$r2").Value2
    $ws.Range("F$r1`:AC$r1").Value = $ws.Range("F$r2`:AC$r2").Value2

    $ws.Range("B$r2").Value        = $tmpB
    $ws.Range("F$r2`:AC$r2").Value = $tmpFAC
}

# Row pairs whose match data was swapped
Swap-Rows 71 72
Swap-Rows 101 102
Swap-Rows 109 110
Swap-Rows 133 134
Swap-Rows 149 150
Swap-Rows 213 214
Swap-Rows 215 216
Swap-Rows 230 231
Swap-Rows 245 246
Swap-Rows 248 249
Swap-Rows 251 252
Swap-Rows 263 265

# ---------------------------------------------------------------------
# New rows appended at the bottom (287-291)
# ---------------------------------------------------------------------

function Set-Row {
    param(
        [int]$r,
        $id,
        $matchId,
        $div,
        $divOrig,
        $date,
        $homeTeam,
        $awayTeam,
        $fthg,
        $ftag,
        $ftr,
        $oddHop,
        $oddDop,
        $oddAop,
        $oddH,
        $oddD,
        $oddA,
        $ah,
        $oddAHH,
        $oddAHA,
        $ahOU,
        $oddAHOver,
        $oddAHUnder,
        $plh,
        $pld,
        $pla,
        $plAhh,
        $plAha,
        $plAhOver,
        $plAhUnder
    )

    # Copy the number/date formats from the last pre-existing row (286)
    # so the new rows carry the same styles (bold/border id cell, date
    # format on the match-date cell) without minting new style entries.
    $ws.Range("A286").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("E286").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$r").Value = $id
    $ws.Range("B$r").Value = $matchId
    $ws.Range("C$r").Value = $div
    $ws.Range("D$r").Value = $divOrig
    $ws.Range("E$r").Value = $date
    $ws.Range("F$r").Value = $homeTeam
    $ws.Range("G$r").Value = $awayTeam
    if ($null -ne $fthg) { $ws.Range("H$r").Value = $fthg }
    if ($null -ne $ftag) { $ws.Range("I$r").Value = $ftag }
    if ($null -ne $ftr)  { $ws.Range("J$r").Value = $ftr }
    $ws.Range("K$r").Value = $oddHop
    $ws.Range("L$r").Value = $oddDop
    $ws.Range("M$r").Value = $oddAop
    $ws.Range("N$r").Value = $oddH
    $ws.Range("O$r").Value = $oddD
    $ws.Range("P$r").Value = $oddA
    $ws.Range("Q$r").Value = $ah
    $ws.Range("R$r").Value = $oddAHH
    $ws.Range("S$r").Value = $oddAHA
    $ws.Range("T$r").Value = $ahOU
    $ws.Range("U$r").Value = $oddAHOver
    $ws.Range("V$r").Value = $oddAHUnder
    $ws.Range("W$r").Value = $plh
    $ws.Range("X$r").Value = $pld
    $ws.Range("Y$r").Value = $pla
    $ws.Range("Z$r").Value = $plAhh
    if ($null -ne $plAha)     { $ws.Range("AA$r").Value = $plAha }
    if ($null -ne $plAhOver)  { $ws.Range("AB$r").Value = $plAhOver }
    if ($null -ne $plAhUnder) { $ws.Range("AC$r").Value = $plAhUnder }
}

# Row 287 - Unam Pumas Women vs Chivas Guadalajara Women (played, FTR=D)
Set-Row 287 285 7645821 "Mexico Liga MX Femenil" "Mexico Liga MX Femenil" `
    45395.625 "Unam Pumas Women" "Chivas Guadalajara Women" `
    1 1 "D" `
    4.5 3.75 1.615 4.333 4 1.571 `
    1 1.775 2.025 3.25 2 1.8 `
    -1 3 -1 0.7749999999999999 -1 -1 0.8

# Row 288 - Leon Women vs Cruz Azul Women (played, FTR=H)
Set-Row 288 286 8089991 "Mexico Liga MX Femenil" "Mexico Liga MX Femenil" `
    45395.92083333333 "Leon Women" "Cruz Azul Women" `
    2 1 "H" `
    1.4 4 7 1.4 4.333 6.5 `
    -1.25 1.825 1.975 3 1.9 1.9 `
    0.3999999999999999 -1 -1 -0.5 0.4875 0 -0

# Row 289 - Tigres UANL Women vs Club America Women (not yet played - no result)
Set-Row 289 287 7645732 "Mexico Liga MX Femenil" "Mexico Liga MX Femenil" `
    45397.91666666666 "Tigres UANL Women" "Club America Women" `
    $null $null $null `
    2.25 3.75 2.5 1.615 4.5 3.6 `
    -0.75 1.825 1.975 3 1.9 1.9 `
    0 0 0 0 0 $null $null

# Row 290 - Santos Laguna Women vs Puebla Women (not yet played - no result)
Set-Row 290 288 7645733 "Mexico Liga MX Femenil" "Mexico Liga MX Femenil" `
    45397.92013888889 "Santos Laguna Women" "Puebla Women" `
    $null $null $null `
    6 5 1.333 4.333 4.5 1.5 `
    1 1.9 1.9 3 1.85 1.95 `
    0 0 0 0 0 $null $null

# Row 291 - Tijuana Women vs Club Necaxa Women (not yet played - no result)
Set-Row 291 289 7645734 "Mexico Liga MX Femenil" "Mexico Liga MX Femenil" `
    45398.00694444445 "Tijuana Women" "Club Necaxa Women" `
    $null $null $null `
    1.166 6.5 10 1.125 9 11 `
    -2.5 1.875 1.925 3.75 1.95 1.85 `
    0 0 0 0 0 $null $null
